$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -ne "小結") {
        $ws.Rows.Item(3).Delete()
        $ws.Range("G2").ClearContents()
        $ws.Range("G7").Value = 999
    }
}
